# Emissao de Apolices e Endossos Page
# Rename the existing sheet and add a new sheet that is a full copy of it
# (same layout, styles, merged cells and logo image), representing the new
# "Cadastro de Parametros de Resseguro" page.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "EnvioApoliceEndossoRE21"

# Duplicate the sheet right after itself so the workbook now has two sheets.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "CadastroParametrosResseguro"

# Update the remembered selection on both sheets.
$ws1.Range("C19").Select()
$ws2.Range("C19").Select()

# Keep the first sheet active/selected.
$ws1.Select()
